# Applies the weekly-refresh edit described by the commit "Fruta / hortaliza, semanal":
# A new week of price data is inserted at row 117, which pushes all the previously
# recorded observations (rows 117-225) down by one row (now rows 118-226). Only the
# columns that vary week to week (D=Fecha, J=Volumen, K=Precio minimo, L=Precio maximo,
# M=Precio promedio ponderado, N=Unidad de comercializacion, P=Precio $/Kg,
# Q=Kg o Unidades) move; the static descriptive columns (A,B,C,E,F,G,H,I,O,R) stay the
# same for every row since they describe the same market/product on every line.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 117
$lastOldRow   = 225
$lastNewRow   = 226

$colsToShift = 4, 10, 11, 12, 13, 14, 16, 17   # D, J, K, L, M, N, P, Q
$staticCols  = 1, 2, 3, 5, 6, 7, 8, 9, 15, 18  # A, B, C, E, F, G, H, I, O, R

# The brand-new last row (226) needs the same static descriptive values (market,
# region, category, etc.) as every other row in this sheet; copy them from the last
# existing row (225) before that row's data columns get shifted away.
foreach ($c in $staticCols) {
    $ws.Cells.Item($lastNewRow, $c).Value2 = $ws.Cells.Item($lastOldRow, $c).Value2
}

# The Fecha (date) column is stored with a date/time number format; make sure the new
# row's D cell carries the same formatting as the rest of the column instead of the
# workbook's plain "General" default.
$ws.Cells.Item($lastNewRow, 4).NumberFormat = $ws.Cells.Item($lastOldRow, 4).NumberFormat

# 1) Snapshot the existing (pre-edit) values for the columns that move, for every row
#    that is about to shift down (rows 117..225).
$snapshot = @{}
for ($r = $firstDataRow; $r -le $lastOldRow; $r++) {
    $rowVals = @{}
    foreach ($c in $colsToShift) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# 2) Shift every snapshotted row down by one: new row (r+1) gets the values that used
#    to live in row r. Go from the bottom up so we never clobber a source row before
#    it has been copied (not strictly required since we already snapshotted, but kept
#    for clarity/robustness).
for ($r = $lastOldRow; $r -ge $firstDataRow; $r--) {
    $destRow = $r + 1
    $rowVals = $snapshot[$r]
    foreach ($c in $colsToShift) {
        $ws.Cells.Item($destRow, $c).Value2 = $rowVals[$c]
    }
}

# 3) Write the brand-new week of data into row 117 (the row freed up by the shift).
$ws.Cells.Item($firstDataRow, 4).Value2  = 44658   # Fecha
$ws.Cells.Item($firstDataRow, 10).Value2 = 100     # Volumen
$ws.Cells.Item($firstDataRow, 11).Value2 = 550     # Precio minimo
$ws.Cells.Item($firstDataRow, 12).Value2 = 600     # Precio maximo
$ws.Cells.Item($firstDataRow, 13).Value2 = 575     # Precio promedio ponderado
$ws.Cells.Item($firstDataRow, 14).Value2 = "$/atado 0,5 a 1 kilo"   # Unidad de comercializacion (unchanged)
$ws.Cells.Item($firstDataRow, 16).Value2 = 575     # Precio $/Kg
$ws.Cells.Item($firstDataRow, 17).Value2 = 1       # Kg o Unidades
